$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 14; this shifts the existing rows 14-44 down to 17-47
# and carries the date-format style of column D along with it.
$ws.Rows("14:16").Insert()

# New row 14: Especial, week of 2023-06-22, Región de O'Higgins
$ws.Cells.Item(14,1).Value = 3
$ws.Cells.Item(14,2).Value = "Femacal de La Calera"
$ws.Cells.Item(14,3).Value = "Coquimbo"
$ws.Cells.Item(14,4).Value = 45099
$ws.Cells.Item(14,5).Value = 5
$ws.Cells.Item(14,6).Value = "Fruta"
$ws.Cells.Item(14,7).Value = 100107
$ws.Cells.Item(14,8).Value = "Otros"
$ws.Cells.Item(14,9).Value = 100107001
$ws.Cells.Item(14,10).Value = "Caqui"
$ws.Cells.Item(14,11).Value = "Mankaki"
$ws.Cells.Item(14,12).Value = "Especial"
$ws.Cells.Item(14,13).Value = 60
$ws.Cells.Item(14,14).Value = 12000
$ws.Cells.Item(14,15).Value = 12000
$ws.Cells.Item(14,16).Value = 12000
$ws.Cells.Item(14,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(14,18).Value = "Región de O'Higgins"
$ws.Cells.Item(14,19).Value = 1200
$ws.Cells.Item(14,20).Value = 10

# New row 15: Primera, week of 2023-06-22, Región de O'Higgins
$ws.Cells.Item(15,1).Value = 3
$ws.Cells.Item(15,2).Value = "Femacal de La Calera"
$ws.Cells.Item(15,3).Value = "Coquimbo"
$ws.Cells.Item(15,4).Value = 45099
$ws.Cells.Item(15,5).Value = 5
$ws.Cells.Item(15,6).Value = "Fruta"
$ws.Cells.Item(15,7).Value = 100107
$ws.Cells.Item(15,8).Value = "Otros"
$ws.Cells.Item(15,9).Value = 100107001
$ws.Cells.Item(15,10).Value = "Caqui"
$ws.Cells.Item(15,11).Value = "Mankaki"
$ws.Cells.Item(15,12).Value = "Primera"
$ws.Cells.Item(15,13).Value = 68
$ws.Cells.Item(15,14).Value = 10000
$ws.Cells.Item(15,15).Value = 10000
$ws.Cells.Item(15,16).Value = 10000
$ws.Cells.Item(15,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(15,18).Value = "Región de O'Higgins"
$ws.Cells.Item(15,19).Value = 1000
$ws.Cells.Item(15,20).Value = 10

# New row 16: Segunda, week of 2023-06-22, Región de O'Higgins
$ws.Cells.Item(16,1).Value = 3
$ws.Cells.Item(16,2).Value = "Femacal de La Calera"
$ws.Cells.Item(16,3).Value = "Coquimbo"
$ws.Cells.Item(16,4).Value = 45099
$ws.Cells.Item(16,5).Value = 5
$ws.Cells.Item(16,6).Value = "Fruta"
$ws.Cells.Item(16,7).Value = 100107
$ws.Cells.Item(16,8).Value = "Otros"
$ws.Cells.Item(16,9).Value = 100107001
$ws.Cells.Item(16,10).Value = "Caqui"
$ws.Cells.Item(16,11).Value = "Mankaki"
$ws.Cells.Item(16,12).Value = "Segunda"
$ws.Cells.Item(16,13).Value = 60
$ws.Cells.Item(16,14).Value = 9000
$ws.Cells.Item(16,15).Value = 9000
$ws.Cells.Item(16,16).Value = 9000
$ws.Cells.Item(16,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(16,18).Value = "Región de O'Higgins"
$ws.Cells.Item(16,19).Value = 900
$ws.Cells.Item(16,20).Value = 10
